$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the summary table
$ws.Range("E83").Value = "Task Type"
$ws.Range("F83").Value = "Count"
$ws.Range("G83").Value = "Total Time"
$ws.Range("H83").Value = "Waiting Time"

# Transportation summary row
$ws.Range("E84").Value = "Transport"
$ws.Range("F84").Formula = '=COUNTIF(A2:A72,"transportation")'
$ws.Range("G84").Formula = '=AVERAGEIF(A2:A72,"transportation",F2:F72)'
$ws.Range("H84").Formula = '=AVERAGEIF(A2:A72,"transportation",H2:H72)'

# Charging summary row
$ws.Range("E85").Value = "Charging"
$ws.Range("F85").Formula = '=COUNTIF(A2:A72,"charging")'
$ws.Range("G85").Formula = '=AVERAGEIF(A2:A72,"charging",F2:F72)'
$ws.Range("H85").Formula = '=AVERAGEIF(A2:A72,"charging",H2:H72)'

# Formatting: header row - bold font, blue fill, border, centered
$headerRange = $ws.Range("E83:H83")
$headerRange.Font.Bold = $true
$headerRange.Interior.ThemeColor = 5
$headerRange.Interior.TintAndShade = 0.599993896298105
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108

# Formatting: data rows - border, centered
$dataRange = $ws.Range("E84:H85")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2
$dataRange.HorizontalAlignment = -4108

$ws.Range("A68").Select()
